# Working version of making input/reactor intervals objects
#
# - components sheet: rename the "glu " input label to "glu", and turn the
#   three previously-empty rows (ace/but/prop intervals) into labelled
#   interval rows (ace, prop, but) with a flag of 1 instead of 0. Add a new
#   trailing (empty) row that keeps the same right-aligned style as the
#   label column so the sheet grows by one row.
# - reactors sheet: the "inputs " column header loses its trailing space.
# - selection/active-sheet view state: components becomes the active sheet
#   and tab, reactors keeps its zoom but its selection/active cell resets.

$wb = $excel.ActiveWorkbook

$components = $wb.Worksheets.Item("components")
$reactors = $wb.Worksheets.Item("reactors")

# --- components sheet ---------------------------------------------------
$components.Range("G2").Value = "glu"

# --- reactors sheet -------------------------------------------------------
$reactors.Range("C1").Value = "inputs"

# --- components sheet (continued) -----------------------------------------
$components.Range("G9").Value = "ace"
$components.Range("G9").HorizontalAlignment = -4152
$components.Range("H9").Value = 1

$components.Range("G10").Value = "prop"
$components.Range("G10").HorizontalAlignment = -4152
$components.Range("H10").Value = 1

$components.Range("G11").Value = "but"
$components.Range("G11").HorizontalAlignment = -4152
$components.Range("H11").Value = 1

# new trailing row, with the same right-aligned style, left otherwise blank
$components.Range("G12").HorizontalAlignment = -4152

# --- view state -----------------------------------------------------------
$reactors.Range("G1").Select()

$components.Activate()
$components.Range("I8").Select()
